$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "76.519.77"
$ws.Range("E2").Value = "  +1.13%  "

# Row 3
$ws.Range("D3").Value = "2.948.85"
$ws.Range("E3").Value = "  +2.48%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'198.78"
$ws.Range("E5").Value = "  +1.83%  "

# Row 6
$ws.Range("D6").Value = "'595.69"
$ws.Range("E6").Value = "  -0.41%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "'0.552"
$ws.Range("E8").Value = "  -0.23%  "

# Row 9
$ws.Range("D9").Value = "'0.203"
$ws.Range("E9").Value = "  +6.22%  "

# Row 10
$ws.Range("D10").Value = "2.949.31"
$ws.Range("E10").Value = "  +2.44%  "

# Row 11
$ws.Range("D11").Value = "'0.442"
$ws.Range("E11").Value = "  +10.07%  "

# Row 12
$ws.Range("E12").Value = "  +0.57%  "

# Row 13
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.492.35"
$ws.Range("E13").Value = "  +1.67%  "

# Row 14
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D14").Value = "'4.90"
$ws.Range("E14").Value = "  -0.18%  "

# Row 15
$ws.Range("D15").Value = "'28.37"
$ws.Range("E15").Value = "  +3.84%  "

# Row 16
$ws.Range("D16").Value = "76.457.53"
$ws.Range("E16").Value = "  +0.98%  "

# Row 17
$ws.Range("D17").Value = "'0.0000191"
$ws.Range("E17").Value = "  +1.08%  "

# Row 18
$ws.Range("D18").Value = "2.955.96"
$ws.Range("E18").Value = "  +1.88%  "

# Row 19
$ws.Range("D19").Value = "'13.51"
$ws.Range("E19").Value = "  +7.84%  "

# Row 20
$ws.Range("D20").Value = "'8.74"
$ws.Range("E20").Value = "  -2.37%  "

# Row 21
$ws.Range("D21").Value = "'377.38"
$ws.Range("E21").Value = "  -0.95%  "

# Row 22
$ws.Range("B22").Value = "SuiNetwork"
$ws.Range("C22").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D22").Value = "'2.29"
$ws.Range("E22").Value = "  -0.08%  "

# Row 23
$ws.Range("B23").Value = "Polkadot"
$ws.Range("C23").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D23").Value = "'4.33"
$ws.Range("E23").Value = "  +4.31%  "

# Row 24
$ws.Range("D24").Value = "'72.24"
$ws.Range("E24").Value = "  +0.70%  "

# Row 25
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.02%  "

# Row 26
$ws.Range("D26").Value = "3.098.51"
$ws.Range("E26").Value = "  +1.82%  "

# Row 27
$ws.Range("D27").Value = "'4.28"
$ws.Range("E27").Value = "  +1.40%  "

# Row 28
$ws.Range("D28").Value = "'9.73"
$ws.Range("E28").Value = "  -0.05%  "

# Row 29
$ws.Range("E29").Value = "  +1.24%  "

# Row 30
$ws.Range("E30").Value = "  +0.22%  "

# Row 31
$ws.Range("D31").Value = "'8.36"
$ws.Range("E31").Value = "  +7.52%  "

# Row 32
$ws.Range("D32").Value = "'1.38"
$ws.Range("E32").Value = "  -1.45%  "

# Row 33
$ws.Range("D33").Value = "'496.93"
$ws.Range("E33").Value = "  -2.17%  "

# Row 34
$ws.Range("E34").Value = "  +0.45%  "

# Row 35
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.08%  "

# Row 36
$ws.Range("D36").Value = "'164.14"
$ws.Range("E36").Value = "  -0.81%  "

# Row 37
$ws.Range("D37").Value = "'20.25"
$ws.Range("E37").Value = "  +0.50%  "

# Row 38
$ws.Range("D38").Value = "'0.393"
$ws.Range("E38").Value = "  +14.13%  "

# Row 39
$ws.Range("E39").Value = "  +21.00%  "

# Row 40
$ws.Range("E40").Value = "  +1.41%  "

# Row 41
$ws.Range("E41").Value = "  -2.84%  "

# Row 42
$ws.Range("E42").Value = "  +0.05%  "

# Row 43
$ws.Range("D43").Value = "'180.33"
$ws.Range("E43").Value = "  -1.56%  "

# Row 44
$ws.Range("D44").Value = "'4.93"
$ws.Range("E44").Value = "  -1.55%  "

# Row 45
$ws.Range("E45").Value = "  -1.24%  "

# Row 46
$ws.Range("E46").Value = "  -0.68%  "

# Row 47
$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").Value = "'1.19"
$ws.Range("E47").Value = "  -2.35%  "

# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "'0.594"
$ws.Range("E48").Value = "  +2.19%  "

# Row 49
$ws.Range("D49").Value = "'2.33"
$ws.Range("E49").Value = "  -0.83%  "

# Row 50
$ws.Range("D50").Value = "'3.88"
$ws.Range("E50").Value = "  +3.11%  "

# Row 51
$ws.Range("E51").Value = "  -0.07%  "
